# Applies the "Finished an excel file with some test data" edit to the
# UncertaintyTest sheet: adds a propagation chain for J (advance ratio),
# reworks a couple of existing formulas/values, and adds rows 15-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UncertaintyTest")

# --- New rows 15-17 labels / header edits, in the order the strings
#     were originally authored (this controls shared-string ordering) --
$ws.Range("A16").Value = "n"
$ws.Range("I1").Value = "dv0/dqu"
$ws.Range("A17").Value = "J"
$ws.Range("A15").Value = "d_p"
$ws.Range("K1").Value = "dJ/dv0"
$ws.Range("L1").Value = "dJ/dn"
$ws.Range("M1").Value = "dJ/dDp"

# --- Row 2: partial-derivative formulas -------------------------------
$ws.Range("B2").Formula = "=B12/2*B14^2"
$ws.Range("H2").Formula = "=(4*SQRT(B5/(B2*B11)+1)*(B9*B10*B14/(2*B11*SQRT(B5/(B2*B11)+1))) + 2*((B11^2)*B14/(B11*SQRT(B5/(B2*B11)+1))-B10*B14)-2*(B9*B10*B2*SQRT(B5/(B2*B11)+1) + 4*(B10^2)*B2*B14*SQRT(B5/(B2*B11)+1) - B10*B5*B14)/(B2*B11*SQRT(B5/(B2*B11)+1)))/(16*(B10^2)*B2*(B5/(B2*B11)+1))"
$ws.Range("I2").Formula = "=(2*B2+(B5/B11))*(B14*B5/(8*B10*(B2^2 + B5*B2/B11)^(3/2)))"
$ws.Range("J2").Formula = "=1+B9/(4*B10) - B5/(4*B10*B2*SQRT(1+B5/(B2*B11)))"
$ws.Range("K2").Formula = "=1/(B16*B15)"
$ws.Range("L2").Formula = "=-B13/((B16^2)*B15)"
$ws.Range("M2").Formula = "=-B13/(B16*(B15^2))"

# --- Row 3: B3 becomes a plain literal instead of "=0.6" --------------
$ws.Range("B3").Value = 0.6

# --- Row 4: B4 becomes a formula ---------------------------------------
$ws.Range("B4").Formula = "=554.67"

# --- Row 5: B5 gains a formula -----------------------------------------
$ws.Range("B5").Formula = "=0.2383*B12*(B16^2)*(B11^4)"

# --- Row 8: B8 value updated -------------------------------------------
$ws.Range("B8").Value = 53.353

# --- Row 12: B12 becomes a plain literal, C12 gains a value ------------
$ws.Range("B12").Value = 0.00183
$ws.Range("C12").Value = 0

# --- Row 13: B13 / C13 gain formulas ------------------------------------
$ws.Range("B13").Formula = "=B14*(1+B9/(4*B10) - B5/(4*B10*B2*SQRT(1+B5/(B2*B11))))"
$ws.Range("C13").Formula = "=SQRT((J2*C14)^2 + (I2*C2)^2 + (H2*C5)^2)"

# --- Row 14: B14 formula reworked ---------------------------------------
$ws.Range("B14").Formula = "=B3*SQRT((B6*B7*B8*B4)/(1+(((B6-1)/2)*(B3^2))))"

# --- Row 15 (new): d_p ---------------------------------------------------
$ws.Range("B15").Formula = "=29.925/12"
$ws.Range("C15").Value = 0.001

# --- Row 16 (new): n ------------------------------------------------------
$ws.Range("B16").Formula = "=5499/60"
$ws.Range("C16").Formula = "=1/60"
$ws.Range("C16").NumberFormat = "General"

# --- Row 17 (new): J --------------------------------------------------------
$ws.Range("B17").Formula = "=B13/(B16*B15)"
$ws.Range("C17").Formula = "=SQRT((K2*C13)^2 + (L2*C16)^2 + (M2*C15)^2)"

# --- Column B width / selection / dimension bookkeeping ---------------------
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Range("V20").Select()
